# Applies the "subimos el último SPA" update to the Semana_7 sheet.
# Updates columns R (uds. Objetivo semana pasada), T (Tendencia Consumo) and
# U (Pedido Final) for the affected article rows, and refreshes the
# Total_Unidades summary cell (C56), which is the sum of column U.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

$updates = @{
    5  = @{ R = 106; T = 0 }
    10 = @{ R = 1;   T = 3 }
    13 = @{ R = 1 }
    23 = @{ R = 2;   T = 1; U = 8 }
    24 = @{ R = 5;   T = 0 }
    25 = @{ R = 9;   T = 0; U = 3 }
    31 = @{ R = 1;   T = 0 }
    33 = @{ R = 4;   T = 7; U = 10 }
    37 = @{ R = 2 }
    39 = @{ R = 4 }
    41 = @{ R = 4;   T = 10; U = 7 }
    44 = @{ R = 1 }
    45 = @{ R = 3;   T = 1; U = 7 }
    46 = @{ R = 2 }
    48 = @{ R = 3 }
    49 = @{ R = 11;  T = 0 }
    50 = @{ R = 2;   T = 0 }
    51 = @{ R = 3 }
    53 = @{ R = 1;   T = 4; U = 2 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $cellRef = "$colLetter$rowNum"
        $ws.Range($cellRef).Value = $cols[$colLetter]
    }
}

# Recalculate the Total_Unidades summary (C56) as the sum of column U
# across the article rows (U3:U53), mirroring the original workbook's
# stored total.
$ws.Range("C56").Value = 124
